$wb = $excel.ActiveWorkbook

# --- 1) "ИсходныеДанные" sheet: add new header columns AX1:BO1 ---
$wsData = $wb.Worksheets.Item("ИсходныеДанные")

$headers = @(
    "IS_Service_type",
    "IS_Service_type_Month",
    "IS_Product_type",
    "IS_Product_type_Month",
    "Pdr_Proj",
    "Pdr_Proj_Month",
    "Proj_Pdr",
    "Proj_Pdr_Month",
    "FN_Month",
    "UHCost_KV1",
    "UMCost_KV1",
    "UHCost_KV2",
    "UMCost_KV2",
    "UHCost_KV3",
    "UMCost_KV3",
    "UHCost_KV4",
    "UMCost_KV4",
    "ISDogName"
)

$startCol = 50  # column AX
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsData.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# --- 2) "Настройки" sheet: add new rows 16-34 describing the new fields ---
$wsSettings = $wb.Worksheets.Item("Настройки")

$settingsRows = @(
    @("IS_Service_type", "Тип сервиса (ИСУ, КИС, ЛИС, ПУ, ..)", "AX", "AY"),
    @("IS_Product_type", "Тип системы (SAP, БК, ЛИМС, MES,…)", "AZ", "BA"),
    @("Pdr_Proj", "Группировка Подразделение+Проект", "BB", "BC"),
    @("Proj_Pdr", "Группировка Проект+Подразделение", "BD", "BE"),
    @("Portfolio", "Портфель проектов", "AV", "AW"),
    @("Personal_email", "Признак отправлять сообщение лично или в общей массе", "AR", $null),
    @("user_email", "Почтовый адрес пользователя", "AS", $null),
    @("boss_email", "Почтовый адрес руководителя данного пользователя", "AT", $null),
    @("Contract", "Доходный договор", "AU", $null),
    @("FN", "Функциональное направление (или подразделение)", "C", "BF"),
    @("UHCost_KV1", "Часовая ставка в 1-м квартале", "BF", $null),
    @("UMCost_KV1", "Месячная ставка в 1-м квартале", "BG", $null),
    @("UHCost_KV2", "Часовая ставка во 2-м квартале", "BH", $null),
    @("UMCost_KV2", "Месячная ставка во 2-м квартале", "BI", $null),
    @("UHCost_KV3", "Часовая ставка в 3-м квартале", "BJ", $null),
    @("UMCost_KV3", "Месячная ставка в 3-м квартале", "BK", $null),
    @("UHCost_KV4", "Часовая ставка в 4-м квартале", "BL", $null),
    @("UMCost_KV4", "Месячная ставка в 4-м квартале", "BM", $null),
    @("ISDogName", "Название ИС из договора", "BO", $null)
)

$startRow = 16
for ($i = 0; $i -lt $settingsRows.Count; $i++) {
    $row = $startRow + $i
    $data = $settingsRows[$i]
    $wsSettings.Cells.Item($row, 1).Value = $data[0]
    $wsSettings.Cells.Item($row, 2).Value = $data[1]
    $wsSettings.Cells.Item($row, 3).Value = $data[2]
    if ($data[3] -ne $null) {
        $wsSettings.Cells.Item($row, 4).Value = $data[3]
    }
}
